$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column F
$ws.Range("F1").Value = "United States"

# Row 2: Lower spending on defense
$ws.Range("B2").Value = 0.157407407407407
$ws.Range("C2").Value = 0.303964757709251
$ws.Range("D2").Value = 0.154761904761905
$ws.Range("E2").Value = 0.284615384615385

# Row 3: Lower spending on retirement pensions
$ws.Range("B3").Value = 0.0416666666666667
$ws.Range("C3").Value = 0.0308370044052863
$ws.Range("D3").Value = 0.00595238095238095
$ws.Range("E3").Value = 0.0307692307692308

# Row 4: Lower spending on healthcare
$ws.Range("B4").Value = 0.037037037037037
$ws.Range("C4").Value = 0.026431718061674
$ws.Range("D4").Value = 0.0119047619047619
$ws.Range("E4").Value = 0.0230769230769231

# Row 5: Lower spending on welfare benefits
$ws.Range("B5").Value = 0.185185185185185
$ws.Range("C5").Value = 0.101321585903084
$ws.Range("D5").Value = 0.0238095238095238
$ws.Range("E5").Value = 0.0923076923076923

# Row 6: Lower spending on education
$ws.Range("B6").Value = 0.0231481481481481
$ws.Range("C6").Value = 0.0220264317180617
$ws.Range("D6").Value = 0.0178571428571429
$ws.Range("E6").Value = 0.00769230769230769

# Row 7: Lower spending on other programs
$ws.Range("B7").Value = 0.231481481481481
$ws.Range("C7").Value = 0.317180616740088
$ws.Range("D7").Value = 0.160714285714286
$ws.Range("E7").Value = 0.107692307692308

# Row 8: Higher taxes on the wealthiest
$ws.Range("B8").Value = 0.532407407407407
$ws.Range("C8").Value = 0.832599118942731
$ws.Range("D8").Value = 0.339285714285714
$ws.Range("E8").Value = 0.876923076923077

# Row 9: Higher corporate income tax rate
$ws.Range("B9").Value = 0.236111111111111
$ws.Range("C9").Value = 0.202643171806167
$ws.Range("D9").Value = 0.148809523809524
$ws.Range("E9").Value = 0.6

# Row 10: Higher personal income tax rates
$ws.Range("B10").Value = 0.0324074074074074
$ws.Range("C10").Value = 0.0704845814977974
$ws.Range("D10").Value = 0.0119047619047619
$ws.Range("E10").Value = 0.261538461538462

# Row 11: Higher public deficit
$ws.Range("B11").Value = 0.0648148148148148
$ws.Range("C11").Value = 0.066079295154185
$ws.Range("D11").Value = 0.0238095238095238
$ws.Range("E11").Value = 0.0923076923076923
